$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4977.6665
$ws.Range("I43").Value = 4773.4
$ws.Range("K43").Value = 4773.4
$ws.Range("M43").Value = -4704.4
$ws.Range("H52").Value = 1540
$ws.Range("I52").Value = 1100
$ws.Range("K52").Value = 3300
$ws.Range("M52").Value = -3140
$ws.Range("H69").Value = 23250
$ws.Range("J69").Value = 23250
$ws.Range("L69").Value = 69750
$ws.Range("N69").Value = -71498
$ws.Range("H72").Value = 23250
$ws.Range("J72").Value = 23250
$ws.Range("L72").Value = 209250
$ws.Range("N72").Value = -217986
$ws.Range("H80").Value = 2402014
$ws.Range("I80").Value = 1783098.9
$ws.Range("J80").Value = 2802488.8
$ws.Range("K80").Value = 5349296.699999999
$ws.Range("L80").Value = 8407466.399999999
$ws.Range("M80").Value = -5348298.699999999
$ws.Range("N80").Value = -8409462.399999999
$ws.Range("H83").Value = 2402014
$ws.Range("I83").Value = 1783098.9
$ws.Range("J83").Value = 2802488.8
$ws.Range("K83").Value = 16047890.1
$ws.Range("L83").Value = 25222399.2
$ws.Range("M83").Value = -16042898.1
$ws.Range("N83").Value = -25232383.2
$ws.Range("H116").Value = 7326.2085
$ws.Range("I116").Value = 6639.75
$ws.Range("K116").Value = 6639.75
$ws.Range("M116").Value = -3197.75
$ws.Range("H132").Value = 4351.1025
$ws.Range("I132").Value = 2731.4075
$ws.Range("J132").Value = 7995.4165
$ws.Range("K132").Value = 8194.2225
$ws.Range("L132").Value = 23986.2495
$ws.Range("M132").Value = -5664.2225
$ws.Range("N132").Value = -29046.2495
$ws.Range("H137").Value = 1971.0667
$ws.Range("I137").Value = 2027.091
$ws.Range("K137").Value = 6081.272999999999
$ws.Range("M137").Value = -3531.272999999999
$ws.Range("H141").Value = 13176346
$ws.Range("I141").Value = 17247076
$ws.Range("K141").Value = 51741228
$ws.Range("M141").Value = -51736048
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2675.611
$ws.Range("I74").Value = 1555.4546
$ws.Range("K74").Value = 1555.4546
$ws.Range("M74").Value = -681.4546
$ws.Range("H77").Value = 2675.611
$ws.Range("I77").Value = 1555.4546
$ws.Range("K77").Value = 7777.273
$ws.Range("M77").Value = -3409.273
$ws.Range("H122").Value = 3711.4119
$ws.Range("I122").Value = 3711.4119
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11134.2357
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -8684.235700000001
$ws.Range("H132").Value = 4662.2915
$ws.Range("I132").Value = 3747.6316
$ws.Range("J132").Value = 8138
$ws.Range("K132").Value = 11242.8948
$ws.Range("L132").Value = 24414
$ws.Range("M132").Value = -8712.8948
$ws.Range("N132").Value = -29474
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 1895
$ws.Range("I26").Value = 1895
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1895
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = $null
$ws.Range("N26").Value = -1608
$ws.Range("H31").Value = 30306152
$ws.Range("I31").Value = 58826668
$ws.Range("J31").Value = 3099.9375
$ws.Range("K31").Value = 58826668
$ws.Range("L31").Value = 3099.9375
$ws.Range("M31").Value = -58826373
$ws.Range("N31").Value = -3689.9375
$ws.Range("H34").Value = 30306152
$ws.Range("I34").Value = 58826668
$ws.Range("J34").Value = 3099.9375
$ws.Range("K34").Value = 58826668
$ws.Range("L34").Value = 3099.9375
$ws.Range("M34").Value = -58826466
$ws.Range("N34").Value = -3503.9375
$ws.Range("H99").Value = 15714.9
$ws.Range("I99").Value = 17292.777
$ws.Range("J99").Value = 1514
$ws.Range("K99").Value = 17292.777
$ws.Range("L99").Value = 1514
$ws.Range("M99").Value = -15794.777
$ws.Range("N99").Value = -4510
$ws.Range("H126").Value = 15714.9
$ws.Range("I126").Value = 17292.777
$ws.Range("J126").Value = 1514
$ws.Range("K126").Value = 51878.33099999999
$ws.Range("L126").Value = 4542
$ws.Range("M126").Value = -49408.33099999999
$ws.Range("N126").Value = -9482
$ws.Range("H132").Value = 3083.4167
$ws.Range("I132").Value = 2972.5557
$ws.Range("J132").Value = 3416
$ws.Range("K132").Value = 8917.667099999999
$ws.Range("L132").Value = 10248
$ws.Range("M132").Value = -6387.667099999999
$ws.Range("N132").Value = -15308
$ws.Range("H134").Value = 2833.1538
$ws.Range("I134").Value = 2735.9583
$ws.Range("J134").Value = 3999.5
$ws.Range("K134").Value = 8207.874899999999
$ws.Range("L134").Value = 11998.5
$ws.Range("M134").Value = -5672.874899999999
$ws.Range("N134").Value = -17068.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 16777.334
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 16777.334
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = $null
$ws.Range("M58").Value = 50332.00199999999
$ws.Range("N58").Value = -50588.00199999999
$ws.Range("H86").Value = 1194.75
$ws.Range("I86").Value = 100
$ws.Range("K86").Value = 300
$ws.Range("M86").Value = 886
$ws.Range("H87").Value = 33333
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").Value = $null
$ws.Range("H89").Value = 1194.75
$ws.Range("I89").Value = 100
$ws.Range("K89").Value = 900
$ws.Range("M89").Value = 5028
$ws.Range("H90").Value = 33333
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").Value = $null
$ws.Range("H113").Value = 1609.6
$ws.Range("I113").Value = 2400
$ws.Range("J113").Value = 1412
$ws.Range("K113").Value = 7200
$ws.Range("L113").Value = 4236
$ws.Range("M113").Value = -5030
$ws.Range("N113").Value = -8576
$ws.Range("H132").Value = 2905
$ws.Range("J132").Value = 2905
$ws.Range("L132").Value = 26145
$ws.Range("N132").Value = -31205
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8359.299999999999
$ws.Range("I70").Value = 2332.3333
$ws.Range("J70").Value = 9422.883
$ws.Range("K70").Value = 2332.3333
$ws.Range("L70").Value = 9422.883
$ws.Range("M70").Value = -2062.3333
$ws.Range("N70").Value = -9962.883
$ws.Range("H73").Value = 8359.299999999999
$ws.Range("I73").Value = 2332.3333
$ws.Range("J73").Value = 9422.883
$ws.Range("K73").Value = 2332.3333
$ws.Range("L73").Value = 9422.883
$ws.Range("M73").Value = -1396.3333
$ws.Range("N73").Value = -11294.883
$ws.Range("H113").Value = 1546049.5
$ws.Range("I113").Value = 3098.7144
$ws.Range("J113").Value = 3706180.5
$ws.Range("K113").Value = 3098.7144
$ws.Range("L113").Value = 3706180.5
$ws.Range("M113").Value = -928.7143999999998
$ws.Range("N113").Value = -3710520.5
$ws.Range("H126").Value = 3562.842
$ws.Range("I126").Value = 3035.7144
$ws.Range("K126").Value = 9107.143199999999
$ws.Range("M126").Value = -6637.143199999999
$ws.Range("H132").Value = 3265
$ws.Range("I132").Value = 3265
$ws.Range("K132").Value = 9795
$ws.Range("M132").Value = -7265
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4590.2856
$ws.Range("I82").Value = 2590.6
$ws.Range("J82").Value = 9589.5
$ws.Range("K82").Value = 2590.6
$ws.Range("L82").Value = 9589.5
$ws.Range("M82").Value = -2229.6
$ws.Range("N82").Value = -10311.5
$ws.Range("H85").Value = 4590.2856
$ws.Range("I85").Value = 2590.6
$ws.Range("J85").Value = 9589.5
$ws.Range("K85").Value = 2590.6
$ws.Range("L85").Value = 9589.5
$ws.Range("M85").Value = -1342.6
$ws.Range("N85").Value = -12085.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 34025
$ws.Range("I74").Value = 33600
$ws.Range("J74").Value = 34237.5
$ws.Range("K74").Value = 33600
$ws.Range("L74").Value = 34237.5
$ws.Range("M74").Value = -32664
$ws.Range("N74").Value = -36109.5
$ws.Range("H77").Value = 34025
$ws.Range("I77").Value = 33600
$ws.Range("J77").Value = 34237.5
$ws.Range("K77").Value = 100800
$ws.Range("L77").Value = 102712.5
$ws.Range("M77").Value = -96120
$ws.Range("N77").Value = -112072.5
$ws.Range("H126").Value = 4159
$ws.Range("J126").Value = 1492
$ws.Range("L126").Value = 4476
$ws.Range("N126").Value = -9416
$ws.Range("H132").Value = 3605.2031
$ws.Range("I132").Value = 3432.2246
$ws.Range("J132").Value = 4170.2666
$ws.Range("K132").Value = 10296.6738
$ws.Range("L132").Value = 12510.7998
$ws.Range("M132").Value = -7766.6738
$ws.Range("N132").Value = -17570.7998
